# Fix missing idAttribute import error by using fully qualified entity names
#
# 1. Rename the "molgenisfieldtypes" sheet to "emx_molgenisfieldtypes".
# 2. Make "emx_molgenisfieldtypes" the active sheet (was "attributes").
# 3. In the "attributes" sheet, the entity / refEntity columns (C and E)
#    referred to the short entity names (molgenisfieldtypes, packagesEMX,
#    entitiesEMX, attributesEMX). Re-point them at the fully qualified
#    "emx_*" names instead.
# 4. Update the remembered selections on the "attributes" and
#    "emx_molgenisfieldtypes" sheets and widen columns C/E on "attributes"
#    to fit the longer fully-qualified names.

$wb = $excel.ActiveWorkbook

$wsAttributes = $wb.Worksheets.Item("attributes")
$wsFieldTypes = $wb.Worksheets.Item("molgenisfieldtypes")

# --- 1. rename sheet ------------------------------------------------------
$wsFieldTypes.Name = "emx_molgenisfieldtypes"

# --- 3. fully-qualify the entity names referenced on "attributes" ---------
$oldToNew = @{
    "molgenisfieldtypes" = "emx_molgenisfieldtypes";
    "packagesEMX"        = "emx_packagesEMX";
    "entitiesEMX"        = "emx_entitiesEMX";
    "attributesEMX"      = "emx_attributesEMX";
}

$usedRange = $wsAttributes.UsedRange
foreach ($row in 1..$usedRange.Rows.Count) {
    foreach ($colLetter in @("C", "E")) {
        $cell = $wsAttributes.Range("$colLetter$row")
        $current = $cell.Value2
        if ($current -and $oldToNew.ContainsKey([string]$current)) {
            $cell.Value = $oldToNew[[string]$current]
        }
    }
}

# --- 4a. widen the now-longer columns --------------------------------------
$wsAttributes.Columns.Item(3).ColumnWidth = 17.65
$wsAttributes.Columns.Item(5).ColumnWidth = 22.65

# --- 2 & 4b. selections + active sheet -------------------------------------
$wsAttributes.Activate()
$wsAttributes.Range("G19").Select()

$wsFieldTypes.Activate()
$wsFieldTypes.Range("H22").Select()
